$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 206
$ws.Range("I9").Value = 162.25
$ws.Range("J9").Value = 249.75
$ws.Range("K9").Value = 162.25
$ws.Range("L9").Value = 249.75
$ws.Range("M9").Value = 6.75
$ws.Range("N9").Value = -587.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4078.8
$ws.Range("J32").Value = 4324.25
$ws.Range("L32").Value = 4324.25
$ws.Range("N32").Value = -4976.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10209.357
$ws.Range("J64").Value = 13888.889
$ws.Range("L64").Value = 13888.889
$ws.Range("N64").Value = -14384.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 10209.357
$ws.Range("J67").Value = 13888.889
$ws.Range("L67").Value = 13888.889
$ws.Range("N67").Value = -15604.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 912
$ws.Range("I132").Value = 862.2162
$ws.Range("J132").Value = 1372.5
$ws.Range("K132").Value = 2586.6486
$ws.Range("L132").Value = 4117.5
$ws.Range("M132").Value = -56.64859999999999
$ws.Range("N132").Value = -9177.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 96142.42999999999
$ws.Range("J136").Value = 96142.42999999999
$ws.Range("L136").Value = 96142.42999999999
$ws.Range("N136").Value = -106342.43

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2102.2126
$ws.Range("I138").Value = 1436
$ws.Range("J138").Value = 2479.7334
$ws.Range("K138").Value = 4308
$ws.Range("L138").Value = 7439.2002
$ws.Range("M138").Value = 832
$ws.Range("N138").Value = -17719.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 59999.5
$ws.Range("J34").Value = 59999.5
$ws.Range("L34").Value = 59999.5
$ws.Range("N34").Value = -60541.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3187.3157
$ws.Range("I45").Value = 1971.0769
$ws.Range("J45").Value = 5822.5
$ws.Range("K45").Value = 1971.0769
$ws.Range("L45").Value = 5822.5
$ws.Range("M45").Value = -1594.0769
$ws.Range("N45").Value = -6576.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7847.2
$ws.Range("J46").Value = 8321.5
$ws.Range("L46").Value = 8321.5
$ws.Range("N46").Value = -8959.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3491.4736
$ws.Range("I74").Value = 2831.72
$ws.Range("J74").Value = 4760.231
$ws.Range("K74").Value = 2831.72
$ws.Range("L74").Value = 4760.231
$ws.Range("M74").Value = -1957.72
$ws.Range("N74").Value = -6508.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3491.4736
$ws.Range("I77").Value = 2831.72
$ws.Range("J77").Value = 4760.231
$ws.Range("K77").Value = 14158.6
$ws.Range("L77").Value = 23801.155
$ws.Range("M77").Value = -9790.599999999999
$ws.Range("N77").Value = -32537.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4894.125
$ws.Range("I132").Value = 4027.45
$ws.Range("K132").Value = 12082.35
$ws.Range("M132").Value = -9552.349999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 108040.664
$ws.Range("J135").Value = 108040.664
$ws.Range("L135").Value = 108040.664
$ws.Range("N135").Value = -118180.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 268.23077
$ws.Range("I11").Value = 109.2
$ws.Range("J11").Value = 798.3333
$ws.Range("K11").Value = 109.2
$ws.Range("L11").Value = 798.3333
$ws.Range("M11").Value = 30.8
$ws.Range("N11").Value = -1078.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2881.9697
$ws.Range("I20").Value = 2088.95
$ws.Range("J20").Value = 4102
$ws.Range("K20").Value = 2088.95
$ws.Range("L20").Value = 4102
$ws.Range("M20").Value = -1841.95
$ws.Range("N20").Value = -4596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4165.6895
$ws.Range("I134").Value = 4250.1787
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 12750.5361
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -10215.5361
$ws.Range("N134").Value = -10470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1833
$ws.Range("I2").Value = 1833
$ws.Range("K2").Value = 1833
$ws.Range("M2").Value = -1720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4478.2964
$ws.Range("I31").Value = 4347.0835
$ws.Range("J31").Value = 4583.2666
$ws.Range("K31").Value = 4347.0835
$ws.Range("L31").Value = 4583.2666
$ws.Range("M31").Value = -4052.0835
$ws.Range("N31").Value = -5173.2666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4478.2964
$ws.Range("I34").Value = 4347.0835
$ws.Range("J34").Value = 4583.2666
$ws.Range("K34").Value = 4347.0835
$ws.Range("L34").Value = 4583.2666
$ws.Range("M34").Value = -4145.0835
$ws.Range("N34").Value = -4987.2666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2113.4285
$ws.Range("I105").Value = 2060.6924
$ws.Range("J105").Value = 2799
$ws.Range("K105").Value = 2060.6924
$ws.Range("L105").Value = 2799
$ws.Range("M105").Value = -313.6923999999999
$ws.Range("N105").Value = -6293

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 939.4
$ws.Range("I107").Value = 829.8
$ws.Range("J107").Value = 1049
$ws.Range("K107").Value = 829.8
$ws.Range("L107").Value = 1049
$ws.Range("M107").Value = 1090.2
$ws.Range("N107").Value = -4889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 224875
$ws.Range("J121").Value = 224875
$ws.Range("L121").Value = 224875
$ws.Range("N121").Value = -227495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3808.95
$ws.Range("J122").Value = 7500
$ws.Range("L122").Value = 22500
$ws.Range("N122").Value = -27400

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5969.385
$ws.Range("I134").Value = 4913.8096
$ws.Range("K134").Value = 14741.4288
$ws.Range("M134").Value = -12206.4288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24020116
$ws.Range("I4").Value = 113493736
$ws.Range("K4").Value = 340481208
$ws.Range("M4").Value = -340481096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1270.1154
$ws.Range("I5").Value = 1190.2778
$ws.Range("K5").Value = 3570.8334
$ws.Range("M5").Value = -3458.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 12694.25
$ws.Range("I23").Value = 225.83333
$ws.Range("K23").Value = 677.49999
$ws.Range("M23").Value = -442.49999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 13166.5
$ws.Range("I126").Value = 9500
$ws.Range("J126").Value = 14999.75
$ws.Range("K126").Value = 28500
$ws.Range("L126").Value = 44999.25
$ws.Range("M126").Value = -23560
$ws.Range("N126").Value = -54879.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 21528730
$ws.Range("J129").Value = 1444914.1
$ws.Range("L129").Value = 4334742.300000001
$ws.Range("N129").Value = -4344742.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1270.1154
$ws.Range("I135").Value = 1190.2778
$ws.Range("K135").Value = 10712.5002
$ws.Range("M135").Value = -8177.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2239.5
$ws.Range("I137").Value = 1737.4
$ws.Range("J137").Value = 4750
$ws.Range("K137").Value = 5212.200000000001
$ws.Range("L137").Value = 14250
$ws.Range("M137").Value = -112.2000000000007
$ws.Range("N137").Value = -24450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2094.476
$ws.Range("I139").Value = 1893.2354
$ws.Range("K139").Value = 5679.706200000001
$ws.Range("M139").Value = -539.7062000000005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11625.223
$ws.Range("I70").Value = 7872.875
$ws.Range("J70").Value = 14627.1
$ws.Range("K70").Value = 7872.875
$ws.Range("L70").Value = 14627.1
$ws.Range("M70").Value = -7602.875
$ws.Range("N70").Value = -15167.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11625.223
$ws.Range("I73").Value = 7872.875
$ws.Range("J73").Value = 14627.1
$ws.Range("K73").Value = 7872.875
$ws.Range("L73").Value = 14627.1
$ws.Range("M73").Value = -6936.875
$ws.Range("N73").Value = -16499.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2252.3713
$ws.Range("I122").Value = 2214.24
$ws.Range("J122").Value = 2347.7
$ws.Range("K122").Value = 6642.719999999999
$ws.Range("L122").Value = 7043.099999999999
$ws.Range("M122").Value = -4192.719999999999
$ws.Range("N122").Value = -11943.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 990.26086
$ws.Range("I16").Value = 814.4211
$ws.Range("K16").Value = 814.4211
$ws.Range("M16").Value = -644.4211

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9505.419
$ws.Range("I46").Value = 3266.6843
$ws.Range("J46").Value = 14444.417
$ws.Range("K46").Value = 3266.6843
$ws.Range("L46").Value = 14444.417
$ws.Range("M46").Value = -3078.6843
$ws.Range("N46").Value = -14820.417

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1358.7
$ws.Range("J68").Value = 1494.6666
$ws.Range("L68").Value = 1494.6666
$ws.Range("N68").Value = -2992.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1358.7
$ws.Range("J71").Value = 1494.6666
$ws.Range("L71").Value = 7473.333000000001
$ws.Range("N71").Value = -14961.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4314
$ws.Range("I122").Value = 3942.5
$ws.Range("J122").Value = 4499.75
$ws.Range("K122").Value = 11827.5
$ws.Range("L122").Value = 13499.25
$ws.Range("M122").Value = -9377.5
$ws.Range("N122").Value = -18399.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11443.429
$ws.Range("I132").Value = 12884.25
$ws.Range("J132").Value = 7841.375
$ws.Range("K132").Value = 38652.75
$ws.Range("L132").Value = 23524.125
$ws.Range("M132").Value = -36122.75
$ws.Range("N132").Value = -28584.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4608.8
$ws.Range("I136").Value = 4390.9473
$ws.Range("J136").Value = 5298.6665
$ws.Range("K136").Value = 13172.8419
$ws.Range("L136").Value = 15895.9995
$ws.Range("M136").Value = -10622.8419
$ws.Range("N136").Value = -20995.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1499.909
$ws.Range("I2").Value = 642.7143
$ws.Range("K2").Value = 642.7143
$ws.Range("M2").Value = -530.7143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9500
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 711.4286
$ws.Range("I107").Value = 414.2
$ws.Range("K107").Value = 1242.6
$ws.Range("M107").Value = 677.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2809.818
$ws.Range("I136").Value = 1974.1666
$ws.Range("K136").Value = 5922.4998
$ws.Range("M136").Value = -3372.4998
